$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 38
$ws.Range("D2").Value = 2
$ws.Range("B5").Value = 0.95
$ws.Range("D5").Value = 0.05
